# Auto-generated Excel COM-interop script
# Applies numeric updates (market price refresh) to the Leve profit sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 826.2222
$ws.Range("I92").Value = 836.7
$ws.Range("J92").Value = 813.125
$ws.Range("K92").Value = 836.7
$ws.Range("L92").Value = 813.125
$ws.Range("M92").Value = 411.3
$ws.Range("N92").Value = -3309.125
$ws.Range("H96").Value = 483.57144
$ws.Range("I96").Value = 542.4545000000001
$ws.Range("J96").Value = 418.8
$ws.Range("K96").Value = 1627.3635
$ws.Range("L96").Value = 1256.4
$ws.Range("M96").Value = -254.3635000000002
$ws.Range("N96").Value = -4002.4
$ws.Range("H100").Value = 1747.0834
$ws.Range("I100").Value = 1531.9286
$ws.Range("J100").Value = 2048.3
$ws.Range("K100").Value = 1531.9286
$ws.Range("L100").Value = 2048.3
$ws.Range("M100").Value = -990.9286
$ws.Range("N100").Value = -3130.3
$ws.Range("H112").Value = 1547.6471
$ws.Range("J112").Value = 1547.6471
$ws.Range("L112").Value = 4642.9413
$ws.Range("N112").Value = -6858.9413
$ws.Range("H129").Value = 660.05554
$ws.Range("I129").Value = 335.125
$ws.Range("J129").Value = 920
$ws.Range("K129").Value = 1005.375
$ws.Range("L129").Value = 2760
$ws.Range("M129").Value = 3994.625
$ws.Range("N129").Value = -12760
$ws.Range("H138").Value = 2269.16
$ws.Range("I138").Value = 1849.421
$ws.Range("J138").Value = 3598.3333
$ws.Range("K138").Value = 5548.263
$ws.Range("L138").Value = 10794.9999
$ws.Range("M138").Value = -408.2629999999999
$ws.Range("N138").Value = -21074.9999
$ws.Range("H141").Value = 3575.375
$ws.Range("J141").Value = 5387.5
$ws.Range("L141").Value = 16162.5
$ws.Range("N141").Value = -26522.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5069.914
$ws.Range("I74").Value = 987.0833
$ws.Range("J74").Value = 13977.909
$ws.Range("K74").Value = 987.0833
$ws.Range("L74").Value = 13977.909
$ws.Range("M74").Value = -113.0833
$ws.Range("N74").Value = -15725.909
$ws.Range("H77").Value = 5069.914
$ws.Range("I77").Value = 987.0833
$ws.Range("J77").Value = 13977.909
$ws.Range("K77").Value = 4935.4165
$ws.Range("L77").Value = 69889.545
$ws.Range("M77").Value = -567.4165000000003
$ws.Range("N77").Value = -78625.545
$ws.Range("H110").Value = 1461.6842
$ws.Range("I110").Value = 1209.2
$ws.Range("J110").Value = 1742.2222
$ws.Range("K110").Value = 1209.2
$ws.Range("L110").Value = 1742.2222
$ws.Range("M110").Value = 835.8
$ws.Range("N110").Value = -5832.2222

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 2608.0715
$ws.Range("I12").Value = 903.3
$ws.Range("J12").Value = 6870
$ws.Range("K12").Value = 903.3
$ws.Range("L12").Value = 6870
$ws.Range("M12").Value = -735.3
$ws.Range("N12").Value = -7206
$ws.Range("H99").Value = 1014.73334
$ws.Range("I99").Value = 799.8889
$ws.Range("J99").Value = 1337
$ws.Range("K99").Value = 799.8889
$ws.Range("L99").Value = 1337
$ws.Range("M99").Value = 698.1111
$ws.Range("N99").Value = -4333
$ws.Range("H102").Value = 20111
$ws.Range("I102").Value = 12638.75
$ws.Range("J102").Value = 50000
$ws.Range("K102").Value = 12638.75
$ws.Range("L102").Value = 50000
$ws.Range("M102").Value = -9393.75
$ws.Range("N102").Value = -56490
$ws.Range("H107").Value = 3864.4285
$ws.Range("I107").Value = 2890.25
$ws.Range("J107").Value = 5163.3335
$ws.Range("K107").Value = 2890.25
$ws.Range("L107").Value = 5163.3335
$ws.Range("M107").Value = -970.25
$ws.Range("N107").Value = -9003.333500000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2030.2
$ws.Range("I107").Value = 2344.6428
$ws.Range("J107").Value = 1296.5
$ws.Range("K107").Value = 2344.6428
$ws.Range("L107").Value = 1296.5
$ws.Range("M107").Value = -424.6428000000001
$ws.Range("N107").Value = -5136.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 12060
$ws.Range("I39").Value = 1900
$ws.Range("J39").Value = 14600
$ws.Range("K39").Value = 5700
$ws.Range("L39").Value = 43800
$ws.Range("M39").Value = -5406
$ws.Range("N39").Value = -44388
$ws.Range("H50").Value = 261
$ws.Range("I50").Value = 308
$ws.Range("J50").Value = 120
$ws.Range("K50").Value = 924
$ws.Range("L50").Value = 360
$ws.Range("M50").Value = -443
$ws.Range("N50").Value = -1322
$ws.Range("H53").Value = 261
$ws.Range("I53").Value = 308
$ws.Range("J53").Value = 120
$ws.Range("K53").Value = 924
$ws.Range("L53").Value = 360
$ws.Range("M53").Value = -443
$ws.Range("N53").Value = -1322
$ws.Range("H131").Value = 2771.59
$ws.Range("I131").Value = 15340
$ws.Range("J131").Value = 1889.5964
$ws.Range("K131").Value = 46020
$ws.Range("L131").Value = 5668.789199999999
$ws.Range("M131").Value = -40980
$ws.Range("N131").Value = -15748.7892
$ws.Range("H137").Value = 21909.299
$ws.Range("I137").Value = 2974.1667
$ws.Range("J137").Value = 26958.666
$ws.Range("K137").Value = 8922.500100000001
$ws.Range("L137").Value = 80875.99800000001
$ws.Range("M137").Value = -3822.500100000001
$ws.Range("N137").Value = -91075.99800000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H99").Value = 11318.167
$ws.Range("I99").Value = 2097.5
$ws.Range("J99").Value = 29759.5
$ws.Range("K99").Value = 2097.5
$ws.Range("L99").Value = 29759.5
$ws.Range("M99").Value = 148.5
$ws.Range("N99").Value = -34251.5
$ws.Range("H113").Value = 1124.125
$ws.Range("J113").Value = 1160
$ws.Range("L113").Value = 1160
$ws.Range("N113").Value = -5500
$ws.Range("H122").Value = 5417.0454
$ws.Range("I122").Value = 6852.769
$ws.Range("J122").Value = 3343.2222
$ws.Range("K122").Value = 20558.307
$ws.Range("L122").Value = 10029.6666
$ws.Range("M122").Value = -18108.307
$ws.Range("N122").Value = -14929.6666

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 10860.3
$ws.Range("I61").Value = 6162.5
$ws.Range("J61").Value = 13992.167
$ws.Range("K61").Value = 6162.5
$ws.Range("L61").Value = 13992.167
$ws.Range("M61").Value = -5870.5
$ws.Range("N61").Value = -14576.167

